$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: a stray '>' that had leaked into the plain-text run right
# after "<del" (the run held "le>" instead of just "le") is removed.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("le>", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "le", 2)

# ------------------------------------------------------------------
# Change 2: "t aprés. Et trois ou" becomes
#           "t aprés. Et <corr>en</corr> trois ou"
# where "<corr>" / "</corr>" are styled like the other red, small-caps
# Courier New markup tags (e.g. "<del>" / "</del>") already present in
# the document, and "en" / " " are plain runs (same convention as the
# plain-styled deleted text already present in the document, e.g. the
# "p" inside "<del>p</del>").
#
# We clone the exact formatting of those existing runs via Copy/Paste
# so the new runs end up with byte-identical <w:rPr> to their sibling
# tag runs elsewhere in the document, then just retarget their text.
# ------------------------------------------------------------------

# Locate the run to split, and the offset right after "Et " (i.e. right
# before "trois ou") where the new markup must be inserted.
$target = $d.Content
$null = $target.Find.Execute("t aprés. Et trois ou", $true, $false, $false, `
                              $false, $false, $true, 1, $false, "", 0)
$splitPos = $target.Start + 12   # Len("t aprés. Et ") == 12

# Locate source runs elsewhere in the document to copy formatting from:
#   "<del>"   -> red/Courier tag-style run (source for "<corr>" / "</corr>")
#   "p"       -> plain run with no explicit color (source for "en" / " ")
#   "</del>"  -> red/Courier tag-style run
$delOpenSrc = $d.Content
$null = $delOpenSrc.Find.Execute("<del>", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
$delOpenRange = $d.Range($delOpenSrc.Start, $delOpenSrc.End)
$plainSrcRange = $d.Range($delOpenSrc.End, $delOpenSrc.End + 1)

$delCloseSrc = $d.Content
$null = $delCloseSrc.Find.Execute("</del>", $true, $false, $false, $false, `
                                   $false, $true, 1, $false, "", 0)
$delCloseRange = $d.Range($delCloseSrc.Start, $delCloseSrc.End)

$pos = $splitPos

# Insert "<corr>" (clone "<del>" formatting, then rename the text)
$delOpenRange.Copy()
$ip = $d.Range($pos, $pos)
$ip.Paste()
$pasted = $d.Range($pos, $pos + 5)          # pasted "<del>" is 5 chars
$pasted.Text = "<corr>"
$pos = $pos + 6                             # "<corr>" is 6 chars

# Insert "en" (clone plain "p" formatting, then rename the text)
$plainSrcRange.Copy()
$ip = $d.Range($pos, $pos)
$ip.Paste()
$pasted = $d.Range($pos, $pos + 1)          # pasted "p" is 1 char
$pasted.Text = "en"
$pos = $pos + 2                             # "en" is 2 chars

# Insert "</corr>" (clone "</del>" formatting, then rename the text)
$delCloseRange.Copy()
$ip = $d.Range($pos, $pos)
$ip.Paste()
$pasted = $d.Range($pos, $pos + 6)          # pasted "</del>" is 6 chars
$pasted.Text = "</corr>"
$pos = $pos + 7                             # "</corr>" is 7 chars

# Insert " " (clone plain "p" formatting again, then rename the text)
$plainSrcRange2 = $d.Range($delOpenSrc.End, $delOpenSrc.End + 1)
$plainSrcRange2.Copy()
$ip = $d.Range($pos, $pos)
$ip.Paste()
$pasted = $d.Range($pos, $pos + 1)          # pasted "p" is 1 char
$pasted.Text = " "
